$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - sheet1
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 262
$ws1.Range("F4").Value = 61
$ws1.Range("F13").Value = 2153
$ws1.Range("F17").Value = 488
$ws1.Range("F19").Value = 78
$ws1.Range("F20").Value = 37
$ws1.Range("F22").Value = 1637
$ws1.Range("F23").Value = 3818
$ws1.Range("F27").Value = 1138
$ws1.Range("F28").Value = 139
$ws1.Range("F29").Value = 2027
$ws1.Range("F32").Value = 79
$ws1.Range("F35").Value = 453
$ws1.Range("F38").Value = 388

# Sheet "演出" (Performance) - sheet2
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 20

# Sheet "全部类型" (All Types) - sheet4
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 262
$ws4.Range("F4").Value = 61
$ws4.Range("F13").Value = 2153
$ws4.Range("F15").Value = 20
$ws4.Range("F18").Value = 488
$ws4.Range("F20").Value = 78
$ws4.Range("F21").Value = 37
$ws4.Range("F23").Value = 1637
$ws4.Range("F24").Value = 3818
$ws4.Range("F28").Value = 1138
$ws4.Range("F29").Value = 139
$ws4.Range("F30").Value = 2027
$ws4.Range("F33").Value = 79
$ws4.Range("F36").Value = 453
$ws4.Range("F39").Value = 388
